$d = $word.ActiveDocument

# Locate the existing bullet that holds the OECD "Nature of Policy Change"
# PDF link -- the new Goodreads bullet is inserted directly after it, as a
# sibling list item (same ListParagraph style / numbering).
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*oecd.org/education/ceri*") {
        $target = $p
    }
}
if ($target -eq $null) {
    throw "Could not locate the OECD hyperlink paragraph to anchor the new bullet"
}

# Insert a new paragraph right after the anchor; it inherits the anchor's
# paragraph formatting (ListParagraph style, ilvl 0 / numId 1).
$anchorRange = $target.Range
$anchorRange.Collapse(0)
$anchorRange.InsertParagraphAfter()

$url = "https://www.goodreads.com/work/quotes/41247321-this-changes-everything-capitalism-vs-the-climate"

# Fill the freshly-created (empty) paragraph with the URL text.
$newPara = $target.Next()
$newPara.Range.Text = $url

# Turn that text into a hyperlink (adds the Hyperlink-styled run + the
# hyperlink relationship), matching the other bulleted links in the doc.
$hlRange = $target.Next().Range
$hlRange.End = $hlRange.End - 1
$d.Hyperlinks.Add($hlRange, $url, "", "", $url) | Out-Null
